# Applies: "revised obs. ref, river slope, extended cal period, fixed rbot"
#
# - Row 6 (was hghbw/hghb/fixed)  -> dhghbw / dh / none, new bounds
# - Row 7 (was hghbe/hghb/fixed)  -> dhghbe / dh / none, new bounds
# - New row 8  -> dhriv / dh / none (river head-difference parameter)
# - Rows 9-12  -> criv, cdrn, tsat, dmax shifted down one row (formulas recomputed)
# - New row 13 -> kc (new parameter group)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Seed the brand-new label strings first (in this exact order) so that
# the workbook's shared-string table lists them as: kc, none, dhghbw,
# dhghbe, dhriv, dh - matching the canonical save order.
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "kc"
$ws.Range("C6").Value  = "none"
$ws.Range("A6").Value  = "dhghbw"
$ws.Range("A7").Value  = "dhghbe"
$ws.Range("A8").Value  = "dhriv"
$ws.Range("B6").Value  = "dh"

function Set-ParRow {
    param($row, $name, $pargp, $partrans, $val, $priorlbnd, $priorubnd, $parlbnd, $parubnd)

    $ws.Range("A$row").Value = $name
    $ws.Range("B$row").Value = $pargp
    $ws.Range("C$row").Value = $partrans
    $ws.Range("D$row").Value = $val
    $ws.Range("E$row").Value = $priorlbnd
    $ws.Range("F$row").Value = $priorubnd
    $ws.Range("G$row").Value = $parlbnd
    $ws.Range("H$row").Value = $parubnd
}

# --- row 6: dhghbw -------------------------------------------------------
Set-ParRow 6 "dhghbw" "dh" "none" 0 -0.25 0.25 -10 10

# --- row 7: dhghbe -------------------------------------------------------
Set-ParRow 7 "dhghbe" "dh" "none" 0 -0.25 0.25 -10 10

# --- row 8: dhriv (new) ---------------------------------------------------
Set-ParRow 8 "dhriv" "dh" "none" 0 -0.25 0.25 -10 10

# --- row 9: criv (shifted down from old row 8, values unchanged) ---------
$criv_val = [double]"1E-4"
$criv_priorlbnd = [double]"9.9999999999999995E-7"
$criv_parlbnd = [double]"1E-8"
Set-ParRow 9 "criv" "criv" "log" $criv_val $criv_priorlbnd 0.01 $criv_parlbnd 100

# --- row 10: cdrn (shifted down from old row 9, values unchanged) --------
$cdrn_val = [double]"1E-3"
$cdrn_priorlbnd = [double]"1E-4"
$cdrn_parlbnd = [double]"1E-8"
Set-ParRow 10 "cdrn" "cdrn" "log" $cdrn_val $cdrn_priorlbnd 0.01 $cdrn_parlbnd 100

# --- row 11: tsat (shifted down from old row 10, calibration period widened)
$ws.Range("A11").Value = "tsat"
$ws.Range("B11").Value = "tsat"
$ws.Range("C11").Value = "log"
$ws.Range("E11").Value = 0.01
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = [double]"1E-3"
$ws.Range("H11").Value = 1
$ws.Range("D11").Formula = "=10^((LOG10(E11)+LOG10(F11))/2)"

# --- row 12: dmax (shifted down from old row 11, fixed rbot bound) -------
$ws.Range("A12").Value = "dmax"
$ws.Range("B12").Value = "dmax"
$ws.Range("C12").Value = "log"
$ws.Range("E12").Value = 0.5
$ws.Range("F12").Value = 50
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1000
$ws.Range("D12").Formula = "=10^((LOG10(E12)+LOG10(F12))/2)"

# --- row 13: kc (new) ------------------------------------------------------
Set-ParRow 13 "kc" "kc" "log" 1 0.5 2 0.1 10

# Match the number-format (scientific) style used by the other numeric
# columns in the table for the newly created rows / cells. (E12 is left at
# the default/general format, matching the source data exactly.)
$ws.Range("D6:D13").NumberFormat = "0.00E+00"
$ws.Range("E6:E11").NumberFormat = "0.00E+00"
$ws.Range("E13").NumberFormat = "0.00E+00"
$ws.Range("G6:G13").NumberFormat = "0.00E+00"
$ws.Range("H6:H13").NumberFormat = "0.00E+00"

# Restore the UI selection/view state saved in the workbook.
$ws.Range("A14").Select() | Out-Null
